# Generate Report for Handoff
#
# The 9c330023-... localization unit (row 3 on each sheet) is queued for a
# new handoff: flip its Status from "Handed back: in sync with en-US" to
# "Ready for handoff" on the Overview sheet (columns B & C) and on each
# language sheet (column B), and stamp the new Latest Handoff Datetime
# (column D) on the per-language sheets.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

# Overview: zh-cn / de-de status columns for the 9c330023 file
$ovw.Range("B3").Value = "Ready for handoff"
$ovw.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: status + new handoff timestamp
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("D3").Value = "2016-02-26 06:33:12"

# de-de sheet: status + new handoff timestamp
$de.Range("B3").Value = "Ready for handoff"
$de.Range("D3").Value = "2016-02-26 06:33:27"
